$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J (rows 2-16)
$data = @(
    @(5, 6),
    @(5, 6),
    @(6, 7),
    @(6, 6),
    @(2, 3),
    @(8, 9),
    @(3, 5),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(2, 2),
    @(2, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
